$wb = $excel.ActiveWorkbook

# --- Sheet2 ("Temple"): add two new sample rows from other runs ------------
$ws2 = $wb.Worksheets.Item("Temple")
$ws2.Select()

# Row 3 - "Arya Stark" / Test-Sample4 / BQ.1.1
$ws2.Cells.Item(3, 1).Value = 44747              # A3 Collection_date
$ws2.Cells.Item(3, 1).NumberFormat = "mm-dd-yy"
$ws2.Cells.Item(3, 2).Value = "Arya Stark"       # B3 Name
$ws2.Cells.Item(3, 5).Value = 4                  # E3 case_id
$ws2.Cells.Item(3, 6).Value = "Test-Sample4"     # F3 SPECIMEN_NUMBER
$ws2.Cells.Item(3, 7).Value = 20.5               # G3 ct value
$ws2.Cells.Item(3, 8).Value = 18                 # H3 age
$ws2.Cells.Item(3, 9).Value = "Female"           # I3 GENDER
$ws2.Cells.Item(3, 14).Value = "Unknown"         # N3 priority
$ws2.Cells.Item(3, 15).Value = "BQ.1.1"          # O3 Actual_lineage

# Row 4 - "Missandei" / Test-Sample5 / XBB.1.5.1
$ws2.Cells.Item(4, 1).Value = 44744              # A4 Collection_date
$ws2.Cells.Item(3, 1).Copy()                     # reuse A3's date style (s=11) instead of minting a new one
$ws2.Cells.Item(4, 1).PasteSpecial(-4122)        # xlPasteFormats
$ws2.Cells.Item(4, 5).Value = 5                  # E4 case_id
$ws2.Cells.Item(4, 6).Value = "Test-Sample5"     # F4 SPECIMEN_NUMBER
$ws2.Cells.Item(4, 7).Value = 26.5               # G4 ct value
$ws2.Cells.Item(4, 8).Value = 25                 # H4 age
$ws2.Cells.Item(4, 9).Value = "Female"           # I4 GENDER
$ws2.Cells.Item(4, 10).Value = "Y"               # J4 breakthrough_case
$ws2.Cells.Item(4, 14).Value = "Breakthrough"    # N4 priority
$ws2.Cells.Item(4, 15).Value = "XBB.1.5.1"       # O4 Actual_lineage
$ws2.Cells.Item(4, 2).Value = "Missandei"        # B4 Name

$ws2.Range("K14").Select()

# --- Sheet1 ("PHL"): only the selection changed -----------------------------
$ws1 = $wb.Worksheets.Item("PHL")
$ws1.Select()
$ws1.Range("I17").Select()

# Temple tab is the one left active in the saved workbook
$ws2.Select()
